$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 1075
$ws.Range("F6").Value = 2773
$ws.Range("F8").Value = 1352
$ws.Range("F9").Value = 948
$ws.Range("F10").Value = 658
$ws.Range("F12").Value = 1237
$ws.Range("F13").Value = 307
$ws.Range("F14").Value = 125
$ws.Range("F15").Value = 773
$ws.Range("F16").Value = 814
$ws.Range("F17").Value = 237
$ws.Range("F18").Value = 573
$ws.Range("F19").Value = 1161
$ws.Range("F21").Value = 691
$ws.Range("F22").Value = 631
$ws.Range("F23").Value = 244
$ws.Range("F24").Value = 337
$ws.Range("F25").Value = 328
$ws.Range("F27").Value = 717
$ws.Range("F28").Value = 8097
$ws.Range("F29").Value = 528
$ws.Range("F30").Value = 528
$ws.Range("F34").Value = 209
$ws.Range("F35").Value = 1679
$ws.Range("F37").Value = 185
$ws.Range("F38").Value = 462
$ws.Range("F41").Value = 205
$ws.Range("F42").Value = 168
$ws.Range("F43").Value = 32
$ws.Range("F44").Value = 84
$ws.Range("F47").Value = 150
$ws.Range("F48").Value = 132
$ws.Range("F50").Value = 88

# Sheet 2: 演出 (Performance)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F8").Value = 75
$ws.Range("F11").Value = 213
$ws.Range("F12").Value = 4424
$ws.Range("F13").Value = 55
$ws.Range("F17").Value = 232

# Sheet 4: 全部类型 (All types)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 1075
$ws.Range("F5").Value = 2773
$ws.Range("F7").Value = 1352
$ws.Range("F8").Value = 948
$ws.Range("F9").Value = 1237
$ws.Range("F12").Value = 125
$ws.Range("F13").Value = 773
$ws.Range("F16").Value = 814
$ws.Range("F17").Value = 237
$ws.Range("F18").Value = 573
$ws.Range("F19").Value = 1161
$ws.Range("F21").Value = 75
$ws.Range("F22").Value = 691
$ws.Range("F23").Value = 631
$ws.Range("F24").Value = 244
$ws.Range("F25").Value = 337
$ws.Range("F26").Value = 328
$ws.Range("F27").Value = 717
$ws.Range("F28").Value = 8098
$ws.Range("F29").Value = 213
$ws.Range("F30").Value = 528
$ws.Range("F31").Value = 528
$ws.Range("F32").Value = 209
$ws.Range("F33").Value = 1679
$ws.Range("F35").Value = 462
$ws.Range("F37").Value = 55
$ws.Range("F38").Value = 55
$ws.Range("F42").Value = 168
$ws.Range("F43").Value = 32
$ws.Range("F44").Value = 84
$ws.Range("F46").Value = 159
$ws.Range("F49").Value = 132
